# Add the "resizer" module row to the module-status table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 3. This shifts every existing row down by two,
# which both makes room for the new module row and preserves the blank-row
# separators that sit between data rows (row 4 and row 9 stay empty, just
# like row 2 and row 7 did before the edit).
$ws.Rows("3:4").Insert()

# Fill in the new module row with its name and description.
$ws.Range("A3").Value = "app.fancy-slider.resizer"
$ws.Range("E3").Value = "In functie de rezolutia maxima si minima definita, acesta calculeaza scalarea sliderul pentru a ocupa tot ecranul."

# Give the new row the same "Good" (green) style used by the other data rows.
$ws.Range("B3:D3").Style = "Good"

# The longer description text means column E needs to be a bit wider.
$ws.Columns("E").ColumnWidth = 100.43

# Move the active selection, matching where the editor ended up.
[void]$ws.Range("C32").Select()
